$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.798601984977722
$ws.Range("B1").Value = 2.210093021392822
$ws.Range("C1").Value = 2.379677057266235
$ws.Range("D1").Value = 3.098793506622314
$ws.Range("E1").Value = 1.455441355705261
